# Re-schedule the MBs for CAN_IDs.
# J1939 IDs re-scheduled at MB16~32 (for receive) and MB32~64 (for transmit).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Standard-ID rows that shift up from MB16..MB20 area into MB8..MB12 (rows 10-14) ---
$ws.Range("B10").Value = "Rx"
$ws.Range("C10").Value = 0
$ws.Range("D10").Value = "STD"
$ws.Range("E10").Value = 0
$ws.Range("F10").Value = "7DF"

$ws.Range("B11").Value = "Rx"
$ws.Range("C11").Value = 0
$ws.Range("D11").Value = "STD"
$ws.Range("E11").Value = 0
$ws.Range("F11").Value = "'7E0"
$ws.Range("F11").NumberFormat = "0.00E+00"

$ws.Range("B12").Value = "Rx"
$ws.Range("C12").Value = 0
$ws.Range("D12").Value = "STD"
$ws.Range("E12").Value = 0
$ws.Range("F12").Value = "'000"

$ws.Range("B13").Value = "Rx"
$ws.Range("C13").Value = 0
$ws.Range("D13").Value = "STD"
$ws.Range("E13").Value = 0
$ws.Range("F13").Value = "7EF"

$ws.Range("B14").Value = "Rx"
$ws.Range("C14").Value = 0
$ws.Range("D14").Value = "STD"
$ws.Range("E14").Value = 0
$ws.Range("F14").Value = "7EE"

# --- Standard-ID transmit pair moves into MB14..MB15 (rows 16-17) ---
$ws.Range("B16").Value = "Tx"
$ws.Range("C16").Value = 0
$ws.Range("D16").Value = "STD"
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = "'5E8"
$ws.Range("F16").NumberFormat = "0.00E+00"

$ws.Range("B17").Value = "Tx"
$ws.Range("C17").Value = 0
$ws.Range("D17").Value = "STD"
$ws.Range("E17").Value = 0
$ws.Range("F17").Value = "'7E8"

# --- J1939 (XTD) receive IDs re-scheduled at MB16~32 (rows 18-30) ---
$ws.Range("B18").Value = "Rx"
$ws.Range("C18").Value = 1
$ws.Range("D18").Value = "XTD"
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = "0x18EA0000"

$ws.Range("B19").Value = "Rx"
$ws.Range("C19").Value = 1
$ws.Range("D19").Value = "XTD"
$ws.Range("E19").Value = 0
$ws.Range("F19").ClearFormats() | Out-Null
$ws.Range("F19").Value = "0x18EB0000"

$ws.Range("B20").Value = "Rx"
$ws.Range("C20").Value = 1
$ws.Range("D20").Value = "XTD"
$ws.Range("E20").Value = 0
$ws.Range("F20").ClearFormats() | Out-Null
$ws.Range("F20").Value = "0x18EC0000"

$ws.Range("B21").Value = "Rx"
$ws.Range("C21").Value = 1
$ws.Range("D21").Value = "XTD"
$ws.Range("E21").Value = 0
$ws.Range("F21").Value = "0x18FEF7D5"

$ws.Range("B22").Value = "Rx"
$ws.Range("C22").Value = 1
$ws.Range("D22").Value = "XTD"
$ws.Range("E22").Value = 0
$ws.Range("F22").Value = "0x18FEF8D5"

$ws.Range("B23").Value = "Rx"
$ws.Range("C23").Value = 1
$ws.Range("D23").Value = "XTD"
$ws.Range("E23").Value = 0
$ws.Range("F23").Value = "0x18FEF9D5"

$ws.Range("B24").Value = "Rx"
$ws.Range("C24").Value = 1
$ws.Range("D24").Value = "XTD"
$ws.Range("E24").Value = 0
$ws.Range("F24").Value = "0x18FEFAD5"

$ws.Range("B25").Value = "Rx"
$ws.Range("C25").Value = 1
$ws.Range("D25").Value = "XTD"
$ws.Range("E25").Value = 0
$ws.Range("F25").Value = "0x18FEFBD5"

$ws.Range("B26").Value = "Rx"
$ws.Range("C26").Value = 1
$ws.Range("D26").Value = "XTD"
$ws.Range("E26").Value = 0
$ws.Range("F26").ClearFormats() | Out-Null
$ws.Range("F26").Value = "0x18FEFCD0"

$ws.Range("B27").Value = "Rx"
$ws.Range("C27").Value = 1
$ws.Range("D27").Value = "XTD"
$ws.Range("E27").Value = 0
$ws.Range("F27").ClearFormats() | Out-Null
$ws.Range("F27").Value = "0x18FEFDD0"

$ws.Range("B28").Value = "Rx"
$ws.Range("C28").Value = 1
$ws.Range("D28").Value = "XTD"
$ws.Range("E28").Value = 0
$ws.Range("F28").Value = "0x18FEFED0"

$ws.Range("B29").Value = "Rx"
$ws.Range("C29").Value = 1
$ws.Range("D29").Value = "XTD"
$ws.Range("E29").Value = 0
$ws.Range("F29").Value = "0x18FEFFD0"

$ws.Range("B30").Value = "Rx"
$ws.Range("C30").Value = 1
$ws.Range("D30").Value = "XTD"
$ws.Range("E30").Value = 0
$ws.Range("F30").Value = "0x08FF10D5"

# --- J1939 (XTD) transmit IDs re-scheduled starting at MB32 (rows 34-35) ---
$ws.Range("B34").Value = "Tx"
$ws.Range("C34").Value = 1
$ws.Range("D34").Value = "XTD"
$ws.Range("E34").Value = 0
$ws.Range("F34").Value = "0x1FFFFFFF"

$ws.Range("B35").Value = "Tx"
$ws.Range("C35").Value = 1
$ws.Range("D35").Value = "XTD"
$ws.Range("E35").Value = 0
$ws.Range("F35").Value = "0x1FFFFFFF"

# --- Rows that no longer carry MB config (now only the MB-name label remains) ---
$ws.Range("B36:F37").ClearContents() | Out-Null
$ws.Range("B50:F63").ClearContents() | Out-Null

# --- Cosmetic: widen column F to fit the longer J1939 IDs, and restore selection ---
$ws.Columns(6).ColumnWidth = 11.29
$ws.Range("G12").Select() | Out-Null
